# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (K instead of Strike#, regen std/mean, s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 0
    10 = 4
    11 = 2
    12 = 0
    13 = 0
    14 = 3
    15 = 0
    16 = 3
    17 = 2
    18 = 4
    19 = 1
    20 = 1
    21 = 0
    22 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
